$d = $word.ActiveDocument

$pairs = @(
    @("486×7=3402", "800×9=7200"),
    @("244×2=488", "197×5=985"),
    @("559×5=2795", "763×2=1526"),
    @("116×4=464", "193×9=1737"),
    @("533×6=3198", "831×7=5817"),
    @("820×2=1640", "867×2=1734"),
    @("362×3=1086", "954×3=2862"),
    @("329×7=2303", "690×3=2070"),
    @("633×3=1899", "589×7=4123"),
    @("161×6=966", "925×5=4625"),
    @("157×8=1256", "266×2=532"),
    @("647×3=1941", "197×9=1773"),
    @("632×7=4424", "114×9=1026"),
    @("979×7=6853", "299×8=2392"),
    @("365×4=1460", "711×2=1422"),
    @("977×2=1954", "220×2=440"),
    @("917×5=4585", "601×6=3606"),
    @("446×6=2676", "172×7=1204"),
    @("786×7=5502", "448×9=4032"),
    @("199×7=1393", "717×3=2151"),
    @("374×2=748", "134×9=1206"),
    @("922×6=5532", "864×7=6048"),
    @("679×7=4753", "388×9=3492"),
    @("378×9=3402", "166×9=1494"),
    @("639×8=5112", "243×9=2187")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false,
                         $true, 1, $false, $new, 2)
}
